# Append 14 new data rows (100-113) to Sheet1, extending the used range
# from A1:K99 to A1:K113, per the "change flowers design and sequence"
# commit: new Training-phase trial rows for participants cn0o2558,
# htd2mqsv, vde2r9iw and e6k1fezv, some carrying a final_total_earnings
# value in column K.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ R=100; A="Anomaly no noise";      B="cn0o2558"; C="Training phase"; D=1; E="['Purple', 'Green']";         F="[['Red', 'Blue'], ['Blue', 'Blue']]";                      G="[None, None]";       H="['8', '10']";     K=$null }
    @{ R=101; A="Anomaly no noise";      B="cn0o2558"; C="Training phase"; D=1; E="['Purple', 'Green']";         F="[['Red', 'Blue'], ['Blue', 'Blue']]";                      G="[None, None]";       H="['8', '10']";     K=$null }
    @{ R=102; A="Anomaly no noise";      B="cn0o2558"; C="Training phase"; D=1; E="['Purple', 'Green']";         F="[['Blue', 'Blue'], ['Blue', 'Red']]";                      G="[None, None]";       H="['10', '8']";     K=$null }
    @{ R=103; A="Anomaly no noise";      B="cn0o2558"; C="Training phase"; D=1; E="['Purple', 'Green']";         F="[['Red', 'Blue'], ['Blue', 'Blue']]";                      G="[None, None]";       H="['8', '10']";     K=$null }
    @{ R=104; A="Anomaly no noise";      B="cn0o2558"; C="Training phase"; D=1; E="['Purple', 'Green']";         F="[['Blue', 'Red'], ['Blue', 'Blue']]";                      G="[None, None]";       H="['8', '10']";     K=$null }
    @{ R=105; A="Anomaly no noise";      B="cn0o2558"; C="Training phase"; D=1; E="['Purple', 'Green']";         F="[['Blue', 'Blue'], ['Red', 'Blue']]";                      G="[None, None]";       H="['10', '8']";     K=$null }
    @{ R=106; A="Anomaly no noise";      B="cn0o2558"; C="Training phase"; D=2; E="['Green', 'Purple']";         F="[['Red', 'Blue'], ['Blue', 'Blue']]";                      G="[None, None]";       H="['8', '10']";     K="1.26" }
    @{ R=107; A="Transmission correct";  B="htd2mqsv"; C="Training phase"; D=1; E="['Purple', 'Green']";         F="[['Red', 'Blue'], ['Blue', 'Yellow']]";                    G="[None, None]";       H="['8', '8']";      K=$null }
    @{ R=108; A="Transmission correct";  B="htd2mqsv"; C="Training phase"; D=1; E="['Purple', 'Green']";         F="[['Red', ''], ['Blue', 'Blue']]";                          G="[None, None]";       H="['2', '10']";     K=$null }
    @{ R=109; A="Transmission correct";  B="htd2mqsv"; C="Training phase"; D=2; E="['Green', 'Purple']";         F="[['Red', 'Red'], ['Blue', 'Blue']]";                       G="[None, None]";       H="['5', '10']";     K=$null }
    @{ R=110; A="Transmission correct";  B="htd2mqsv"; C="Training phase"; D=3; E="['Purple', 'Green', 'Green']"; F="[['Red', ''], ['Yellow', 'Yellow'], ['Blue', 'Blue']]";   G="[None, None, None]"; H="['2', '6', '10']"; K=$null }
    @{ R=111; A="Transmission correct";  B="htd2mqsv"; C="Training phase"; D=4; E="['Green', 'Purple', 'Purple']"; F="[['Red', ''], ['Blue', ''], ['Blue', 'Blue']]";          G="[None, None, None]"; H="['2', '5', '10']"; K="0.78" }
    @{ R=112; A="Anomaly no noise";      B="vde2r9iw"; C="Training phase"; D=1; E="['Purple', 'Green']";         F="[['Red', ''], ['Blue', '']]";                              G="[None, None]";       H="['2', '5']";      K="0.07" }
    @{ R=113; A="Anomaly no noise";      B="e6k1fezv"; C="Training phase"; D=1; E="['Purple', 'Green']";         F="[['Red', ''], ['Blue', '']]";                              G="[None, None]";       H="['2', '5']";      K="0.07" }
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    if ($null -ne $row.K) {
        # Leading apostrophe forces Excel to store this numeric-looking
        # value as text (matches the source inlineStr cell, not a number).
        $ws.Cells.Item($r, 11).Value = "'" + $row.K
        # Drop the "quote-prefix" text style Excel auto-applies above so
        # the cell keeps the workbook's plain default formatting.
        $ws.Cells.Item($r, 11).Style = "Normal"
    }
}
